$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. CONCERN.FINE (column I) fine-tuning: shift several bands down by one level ---

# Rows 134-145: I 5 -> 4
for ($r = 134; $r -le 145; $r++) {
    $ws.Cells.Item($r, 9).Value = 4
}

# Rows 275-281: I 4 -> 3
for ($r = 275; $r -le 281; $r++) {
    $ws.Cells.Item($r, 9).Value = 3
}

# Rows 416-435: I 3 -> 2
for ($r = 416; $r -le 435; $r++) {
    $ws.Cells.Item($r, 9).Value = 2
}

# Rows 551-584: I 2 -> 1
for ($r = 551; $r -le 584; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# --- 2. CONCERN.COARSE (column K) recategorised MID -> LOW for rows 453-488 ---
for ($r = 453; $r -le 488; $r++) {
    $ws.Cells.Item($r, 11).Value = "LOW"
}

# --- 3. INV.C (column G) recomputed float precision refresh for rows 654-663 ---
$ws.Cells.Item(654, 7).Value = 0.6800000000000001
$ws.Cells.Item(655, 7).Value = 0.6800000000000001
$ws.Cells.Item(656, 7).Value = 0.69
$ws.Cells.Item(657, 7).Value = 0.69
$ws.Cells.Item(658, 7).Value = 0.69
$ws.Cells.Item(659, 7).Value = 0.69
$ws.Cells.Item(660, 7).Value = 0.69
$ws.Cells.Item(661, 7).Value = 0.69
$ws.Cells.Item(662, 7).Value = 0.69
$ws.Cells.Item(663, 7).Value = 0.69

# --- 4. Drop trailing low-volume districts (rows 700-737) entirely ---
$ws.Rows("700:737").Delete()
